$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match row appended at the bottom of the results table (row 98),
# following the same column layout as the existing rows (2..97):
# Indice, pais, torneio, temporada, data_partida, home, home_ft_gols,
# away, away_ft_gols, home_opening_odds, home_opening_data_hora,
# home_closing_odds, home_closing_data_hora, draw_opening_odds,
# draw_opening_data_hora, draw_closing_odds, draw_closing_data_hora,
# away_opening_odds, away_opening_data_hora, away_closing_odds,
# away_closing_data_hora, url_partida.

# Column A (Indice) reuses the bold/centered/bordered style already used by
# A2:A97 -- copy the format from A97 so the existing style is shared rather
# than duplicated, then overwrite with the real value.
$ws.Range("A97").Copy()
$ws.Range("A98").PasteSpecial(-4122)
$ws.Range("A98").Value = 97

$ws.Range("B98").Value = "ecuador"
$ws.Range("C98").Value = "liga-pro"

# "2023" must stay text (like D2:D97), not become the number 2023 -- the
# leading apostrophe forces a text literal, same as typing it in the UI.
$ws.Range("D98").Value = "'2023"
$ws.Range("D98").ClearFormats()

# Column E (data_partida) reuses the date/time number format already used by
# E2:E97 -- copy the format from E97, then overwrite with the real value.
$ws.Range("E97").Copy()
$ws.Range("E98").PasteSpecial(-4122)
$ws.Range("E98").Value = 45241.04166666666

$ws.Range("F98").Value = "Orense"
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = "Mushuc Runa"
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 1.85
$ws.Range("K98").Value = "05/11/2023 21:12"
$ws.Range("L98").Value = 1.71
$ws.Range("M98").Value = "11/11/2023 00:54"
$ws.Range("N98").Value = 3.48
$ws.Range("O98").Value = "05/11/2023 21:12"
$ws.Range("P98").Value = 3.7
$ws.Range("Q98").Value = "11/11/2023 00:54"
$ws.Range("R98").Value = 4.36
$ws.Range("S98").Value = "05/11/2023 21:12"
$ws.Range("T98").Value = 5.13
$ws.Range("U98").Value = "11/11/2023 00:55"
$ws.Range("V98").Value = "https://www.betexplorer.com/football/ecuador/liga-pro/orense-mushuc-runa/zHHrEcJk/"
